$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text would otherwise be auto-coerced to a number by Excel
# (losing exact formatting such as trailing zeros) get a text NumberFormat first.
$textCells = @("D5","D6","D7","D11","D12","D13","D14","D16","D18","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D34","D35","D37","D41","D42","D43","D44","D45","D46","D47","D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell.
$ws.Range("D2").Value = "66.494.55"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "3.341.93"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "190.40"
$ws.Range("E5").Value = "  +4.73%  "
$ws.Range("D6").Value = "567.91"
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "3.334.48"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("D11").Value = "0.592"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("D12").Value = "48.09"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").Value = "0.0000274"
$ws.Range("E13").Value = "  +3.88%  "
$ws.Range("D14").Value = "8.74"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "3.872.28"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "609.64"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "66.512.00"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "18.15"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").Value = "3.333.98"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "11.22"
$ws.Range("E21").Value = "  -1.72%  "
$ws.Range("D22").Value = "0.919"
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("D23").Value = "18.91"
$ws.Range("E23").Value = "  +12.14%  "
$ws.Range("E24").Value = "  +3.89%  "
$ws.Range("D25").Value = "101.31"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("D26").Value = "4.05"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "2.78"
$ws.Range("E27").Value = "  +4.29%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "5.98"
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("D29").Value = "9.82"
$ws.Range("E29").Value = "  +5.83%  "
$ws.Range("D30").Value = "8.76"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").Value = "30.66"
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("D32").Value = "6.86"
$ws.Range("E32").Value = "  +9.91%  "
$ws.Range("E33").Value = "  +7.56%  "
$ws.Range("D34").Value = "571.05"
$ws.Range("E34").Value = "  +6.07%  "
$ws.Range("D35").Value = "11.17"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("E36").Value = "  +2.19%  "
$ws.Range("D37").Value = "57.37"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").Value = "3.720.48"
$ws.Range("E38").Value = "  -3.26%  "
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").Value = "0.0₃0735"
$ws.Range("E40").Value = "  +3.50%  "
$ws.Range("D41").Value = "34.28"
$ws.Range("E41").Value = "  +7.23%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "3.33"
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.132"
$ws.Range("E43").Value = "  +5.41%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "2.73"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("B45").Value = "CoreDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D45").Value = "3.43"
$ws.Range("E45").Value = "  +9.92%  "
$ws.Range("D46").Value = "0.347"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").Value = "0.0430"
$ws.Range("E47").Value = "  +4.41%  "
$ws.Range("D48").Value = "3.30"
$ws.Range("E48").Value = "  +5.89%  "
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("E51").Value = "  -0.09%  "
